# Apply the "early stopping / new hyperparam search" edit:
#  1. Insert three new columns (F, G, H) for PGD-training hyperparameters.
#  2. Rename the old "adv_train_flag" header to "training_type".
#  3. Label the three new header cells.
#  4. Re-blank the (now-shifted) body cells in F:H for the pre-existing rows
#     (they have no PGD hyperparameters recorded).
#  5. Append three new data rows (25-27) for the new PGD hyperparameter-search
#     runs, including values for the new F/G/H columns, and leave their
#     Distances/Norms (S/T) columns blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 3 columns before the old "activation" column (F) ---------
$ws.Columns("F:H").Insert()

# --- 2 & 3. Header row text ------------------------------------------------
$ws.Range("E1").Value = "training_type"
$ws.Range("F1").Value = "pgd_train_eps"
$ws.Range("G1").Value = "pgd_train_eta"
$ws.Range("H1").Value = "pgd_train_num_iter"

# --- 4. Blank out F:H for the existing data rows (2-24) --------------------
# These runs were recorded before PGD-training hyperparameters existed, so
# the new columns are empty for them (same as the other "not applicable"
# cells already in the sheet).
$ws.Range("F2:H24").Style = "Normal"

# --- 5. Append the three new rows (25-27) -----------------------------------
$ws.Range("A25").Value = 23
$ws.Range("B25").Value = 2
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 0.003
$ws.Range("E25").Value = "PGD"
$ws.Range("F25").Value = 0.1
$ws.Range("G25").Value = 0.1
$ws.Range("H25").Value = 3
$ws.Range("I25").Value = "<function relu at 0x118b969d8>"
$ws.Range("J25").Value = 0.8787999749183655
$ws.Range("K25").Value = 0.01040000002831221
$ws.Range("L25").Value = 0.001300000003539026
$ws.Range("M25").Value = 0.4310351014137268
$ws.Range("N25").Value = 8.296195983886719
$ws.Range("O25").Value = 0.01040000002831221
$ws.Range("P25").Value = "logs/results_278.log"
$ws.Range("Q25").Value = "weights/model_278.ckpt"
$ws.Range("R25").Value = "tb/278/robust"

$ws.Range("A26").Value = 24
$ws.Range("B26").Value = 2
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 0.003
$ws.Range("E26").Value = "PGD"
$ws.Range("F26").Value = 0.1
$ws.Range("G26").Value = 0.1
$ws.Range("H26").Value = 3
$ws.Range("I26").Value = "<function relu at 0x121b0f9d8>"
$ws.Range("J26").Value = 0.8522999882698059
$ws.Range("K26").Value = 0.002899999963119626
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 0.5192863941192627
$ws.Range("N26").Value = 9.538139343261719
$ws.Range("O26").Value = 0.002899999963119626
$ws.Range("P26").Value = "logs/results_279.log"
$ws.Range("Q26").Value = "weights/model_279.ckpt"
$ws.Range("R26").Value = "tb/279/robust"

$ws.Range("A27").Value = 25
$ws.Range("B27").Value = 2
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 0.003
$ws.Range("E27").Value = "PGD"
$ws.Range("F27").Value = 0.1
$ws.Range("G27").Value = 0.1
$ws.Range("H27").Value = 3
$ws.Range("I27").Value = "<function relu at 0x121b0f9d8>"
$ws.Range("J27").Value = 0.8676999807357788
$ws.Range("K27").Value = 0.004999999888241291
$ws.Range("L27").Value = [double]"9.999999747378752e-05"
$ws.Range("M27").Value = 0.4790646433830261
$ws.Range("N27").Value = 9.238405227661133
$ws.Range("O27").Value = 0.004999999888241291
$ws.Range("P27").Value = "logs/results_279.log"
$ws.Range("Q27").Value = "weights/model_279.ckpt"
$ws.Range("R27").Value = "tb/279/robust"

# S25:T27 (Distances / Norms) are left blank for the new rows - mark them as
# present-but-empty cells, matching the rest of the sheet's "N/A" styling.
$ws.Range("S25:T27").Style = "Normal"
